$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: extend merged ranges from Q to S for rows 1-5
$ws.Range("A1:Q1").UnMerge()
$ws.Range("A2:Q2").UnMerge()
$ws.Range("A3:Q3").UnMerge()
$ws.Range("A4:Q4").UnMerge()
$ws.Range("A5:Q5").UnMerge()

$ws.Range("Q1").Copy()
$ws.Range("R1:S1").PasteSpecial(-4122)
$ws.Range("Q2").Copy()
$ws.Range("R2:S2").PasteSpecial(-4122)
$ws.Range("Q3").Copy()
$ws.Range("R3:S3").PasteSpecial(-4122)
$ws.Range("Q4").Copy()
$ws.Range("R4:S4").PasteSpecial(-4122)
$ws.Range("Q5").Copy()
$ws.Range("R5:S5").PasteSpecial(-4122)

$ws.Range("A1:S1").Merge()
$ws.Range("A2:S2").Merge()
$ws.Range("A3:S3").Merge()
$ws.Range("A4:S4").Merge()
$ws.Range("A5:S5").Merge()

# Step 2: Add headers to R6 and S6, copying format from Q6
$ws.Range("Q6").Copy()
$ws.Range("R6:S6").PasteSpecial(-4122)
$ws.Range("R6").Value = "{col:diecisiete}"
$ws.Range("S6").Value = "{col:dieciocho}"

Write-Host "done"
